$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "TempSheet"
